# Apply the "Liste des taches" update:
#  - Fill in missing "Date de debut" / "Date de fin" values for several tasks
#    (reusing the existing date-formatted style from neighbouring cells)
#  - Move the "Termine" (X) marker for two tasks from "A faire" to "Termine"
#  - Update the saved selection on the first sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 9 / 10: fill in "Date de fin" (F) to match the existing "Date de debut" (G) style ---
$ws.Range("G9").Copy()
$ws.Range("F9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F9").Value = 42512

$ws.Range("G10").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F10").Value = 42512

# --- Row 13: fill in the missing end date (G13), copying the date style from F13 ---
$ws.Range("F13").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("G13").Value = 42516

# --- Row 14: fill in both dates, copying the date style from F13 ---
$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 42516

$ws.Range("F13").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("G14").Value = 42516

# --- Row 16: move the "X" marker from "A faire" (B) to "Termine" (D), add both dates ---
$ws.Range("B16").ClearContents()
$ws.Range("D16").Value = "X"

$ws.Range("F15").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value = 42517

$ws.Range("G15").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 42520

# --- Row 17: same change as row 16 ---
$ws.Range("B17").ClearContents()
$ws.Range("D17").Value = "X"

$ws.Range("F15").Copy()
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F17").Value = 42517

$ws.Range("G15").Copy()
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("G17").Value = 42520

# --- Row 21: fill in both dates ---
$ws.Range("F15").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F21").Value = 42514

$ws.Range("G15").Copy()
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("G21").Value = 42514

# --- Row 28: fill in the missing "Date de fin" (F), matching the "Date de debut" (G) style ---
$ws.Range("G28").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 42512

$excel.CutCopyMode = 0

# --- Update the active view/selection on the sheet ---
$ws.Activate()
$ws.Range("G21").Select()

# --- Window size (best-effort; headless runtime has no real window chrome) ---
$excel.ActiveWindow.Width = 14835
$excel.ActiveWindow.Height = 6240
